$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2017-02-09 09:20:16"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2017-02-09 09:19:54"
$wsZhCn.Range("L2").Value = "2017-02-09 09:20:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2017-02-09 09:20:16"
$wsDeDe.Range("L2").Value = "2017-02-09 09:21:25"
